$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '64.154.56'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -5.55%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.293.10'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -7.61%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '179.12'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -13.39%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '516.26'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -8.22%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.592'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -2.84%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.288.26'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -7.56%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.617'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -8.44%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '57.34'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -5.21%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.131'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -10.30%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000255'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -8.21%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.05'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -10.97%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.806.81'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -8.28%  '
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -6.05%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.277.75'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -8.44%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '63.805.41'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -5.78%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '17.13'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -8.64%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.79'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -11.04%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.945'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -10.34%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '370.98'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -7.23%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '79.89'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -4.86%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.63'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -11.89%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '10.84'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -12.93%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -1.71%  '
$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '6.02'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -1.74%  '
$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.64'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -7.59%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '11.22'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -9.08%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.28'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -9.39%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '28.43'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -9.12%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '643.41'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -2.63%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.67'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -13.09%  '
$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '11.05'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -7.90%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '59.16'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -6.21%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.103'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -8.08%  '
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.21%  '
$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '35.73'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -12.54%  '
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.375'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -7.63%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.994'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -0.51%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.121'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -8.45%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.846.66'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -9.72%  '
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0₃0653'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -12.39%  '
$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.65'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -18.93%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.60'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -6.75%  '
$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.32'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -13.29%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0380'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -6.72%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.74'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +2.38%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.123'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -4.77%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.92'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -3.55%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '132.49'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -4.25%  '
